# Insert 9 new daily rows (2019-11-18 .. 2019-11-28) into the historical
# price table, right before the existing 2019-11-29 row (currently row 322),
# shifting all subsequent rows down by 9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room: insert 9 blank rows starting at row 322.
$ws.Range("A322:A330").EntireRow.Insert()

$newRows = @(
    @(1574035200, "2019-11-18", "0202", "RGTECH", 0.29,  0.3,   0.285, 0.295, 2064000),
    @(1574121600, "2019-11-19", "0202", "RGTECH", 0.295, 0.31,  0.295, 0.305, 4405900),
    @(1574208000, "2019-11-20", "0202", "RGTECH", 0.31,  0.315, 0.305, 0.31,  5015500),
    @(1574294400, "2019-11-21", "0202", "RGTECH", 0.31,  0.315, 0.3,   0.31,  1272000),
    @(1574380800, "2019-11-22", "0202", "RGTECH", 0.31,  0.31,  0.295, 0.305, 662700),
    @(1574640000, "2019-11-25", "0202", "RGTECH", 0.305, 0.31,  0.29,  0.3,   3093400),
    @(1574726400, "2019-11-26", "0202", "RGTECH", 0.3,   0.335, 0.3,   0.33,  23573200),
    @(1574812800, "2019-11-27", "0202", "RGTECH", 0.33,  0.335, 0.32,  0.325, 4942300),
    @(1574899200, "2019-11-28", "0202", "RGTECH", 0.33,  0.335, 0.32,  0.33,  3641300)
)

$startRow = 322
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    # Leading apostrophe forces these to stay plain text (matches the
    # original file's inlineStr date/id columns) instead of being
    # auto-converted to a date serial number / stripped of leading zero.
    $ws.Cells.Item($r, 2).Value = "'" + $row[1]
    $ws.Cells.Item($r, 3).Value = "'" + $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
}
